# Fruta / hortaliza, semanal
# Insert a new weekly price record as row 19 (Castle Brite / Especial,
# Damasco, Mercado Mayorista Lo Valledor de Santiago). All subsequent
# rows (previously 19-48) shift down one position to 20-49.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 19, pushing existing rows 19-48 down to 20-49.
$ws.Rows(19).Insert()

# Populate the new row 19 with the new data point.
$ws.Range("A19").Value = 6
$ws.Range("B19").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C19").Value = "Metropolitana"
$ws.Range("D19").Value = 44526
$ws.Range("E19").Value = 13
$ws.Range("F19").Value = "Fruta"
$ws.Range("G19").Value = 100103
$ws.Range("H19").Value = "Frutos de hueso (carozo)"
$ws.Range("I19").Value = 100103003
$ws.Range("J19").Value = "Damasco"
$ws.Range("K19").Value = "Castle Brite"
$ws.Range("L19").Value = "Especial"
$ws.Range("M19").Value = 240
$ws.Range("N19").Value = 25000
$ws.Range("O19").Value = 25000
$ws.Range("P19").Value = 25000
$ws.Range("Q19").Value = "`$/caja 16 kilos"
$ws.Range("R19").Value = "Región de O'Higgins"
$ws.Range("S19").Value = 1562
$ws.Range("T19").Value = 16
